# Template update: reorder the K1:M1 header columns (Lateralitaet / Herdlaesion /
# Zweitlaesion), tidy up the resulting column widths, and refresh the view so the
# sheet opens scrolled to A1 with column N selected instead of cell AG1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the K1:M1 header cells -------------------------------------
# Move the "Lateralität" column (was M) so it becomes the new K column;
# "Herdläsion" (was K) and "Zweitläsion" (was L) each slide one column right.
$ws.Range("M1").EntireColumn.Cut()
$ws.Range("K1").EntireColumn.Insert()

# --- Column widths: K & L now share the same custom width ---------------
$ws.Range("K1:L1").EntireColumn.ColumnWidth = 9.166666666666666

# --- Update the view: scroll back to A1 and select column N -------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N1:N1048576").Select()
